$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 13, shifting existing rows 13-23 down to 14-24.
$ws.Rows.Item(13).Insert()

# The inserted row copies column A's (bold) formatting into A13. Re-apply the normal
# value/highlight formatting (columns B/C) to B13/C13, then clear A13 so it stays empty,
# matching how the other "value" rows (no label in column A) look, e.g. row 24.
$ws.Range("B14:C14").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("A13").Clear()

# Row 10 (Objetivos:) - B/C previously held the teacher's name by mistake; replace with the
# actual objectives text.
$ws.Range("B10").Value = "Apresentar conceitos fundamentais de Administração da Produção."
$ws.Range("C10").Value = "Apresentar conceitos fundamentais de Administração da Produção."

# Row 13 (new, blank row under "Docentes responsáveis:") - fill in the teacher's name.
$ws.Range("B13").Value = "5840560 - Marco Antonio Carvalho Pereira"
$ws.Range("C13").Value = "5840560 - Marco Antonio Carvalho Pereira"

# Row 14 (Programa resumido:) - replace placeholder "Semestral" with the real short syllabus.
$ws.Range("B14").Value = "1. Fundamentos da Gestão de Produção2. Visão estratégica da Produção.3. Projeto em Gestão da Produção.4. Planejamento e Controle da Produção"
$ws.Range("C14").Value = "1. Fundamentos da Gestão de Produção2. Visão estratégica da Produção.3. Projeto em Gestão da Produção.4. Planejamento e Controle da Produção"

# Row 16 (Programa:) - replace placeholder date with the real syllabus text.
$ws.Range("B16").Value = "1 - Fundamentos da gestão de produção: modelo de transformação: inputs, processo de transformação e outputs. Tipos de Processo de Produção2 - Visão estratégica de produção: Papel da função produção. Objetivos de Desempenho. Estratégias de Produção. Ciclo de Vida Produto/Serviço.3 – Projeto em Gestão da Produção: Tipos de Processos. Projeto de Produtos e Serviços. Projeto de Rede de Operações Produtivas. Arranjo Físico.4 - Planejamento e Controle da Produção: Material Requirement Planning (MRP), Manufacturing Resources Planning (MPRII), Enterprise Planning (ERP). Produção Enxuta. Kanban. Just in Time."
$ws.Range("C16").Value = "1 - Fundamentos da gestão de produção: modelo de transformação: inputs, processo de transformação e outputs. Tipos de Processo de Produção2 - Visão estratégica de produção: Papel da função produção. Objetivos de Desempenho. Estratégias de Produção. Ciclo de Vida Produto/Serviço.3 – Projeto em Gestão da Produção: Tipos de Processos. Projeto de Produtos e Serviços. Projeto de Rede de Operações Produtivas. Arranjo Físico.4 - Planejamento e Controle da Produção: Material Requirement Planning (MRP), Manufacturing Resources Planning (MPRII), Enterprise Planning (ERP). Produção Enxuta. Kanban. Just in Time."

# Row 19 (Método:) - previously held the teacher's name by mistake; replace with the real method text.
$ws.Range("B19").Value = "O desenvolvimento da disciplina será baseado em leituras, aula expositiva, discussão e resolução de estudos de caso."
$ws.Range("C19").Value = "O desenvolvimento da disciplina será baseado em leituras, aula expositiva, discussão e resolução de estudos de caso."

# Row 20 (Critério:) - replace with the real criterion text.
$ws.Range("B20").Value = "Provas e Trabalhos"
$ws.Range("C20").Value = "Provas e Trabalhos"

# Row 21 (Norma de recuperação:) - replace with the real recovery-norm text.
$ws.Range("B21").Value = "Prova única com nota maior ou igual a 5,0 (cinco)."
$ws.Range("C21").Value = "Prova única com nota maior ou igual a 5,0 (cinco)."

# Row 22 (Bibliografia:) - replace with the real bibliography text.
$ws.Range("B22").Value = "SLACK, N. et al. Administração da Produção. 3 ed. São Paulo: Atlas, 2009.CHASE, R. B. E JACOBS, F.R. Administração da Produção e de Operações. 1 ed. Porto Alegre. Bookman. 2009.CORREA, H.L.; CORREA, C.A. Administração da Produção e Operações. 2 ed. São Paulo. Atlas. 2006"
$ws.Range("C22").Value = "SLACK, N. et al. Administração da Produção. 3 ed. São Paulo: Atlas, 2009.CHASE, R. B. E JACOBS, F.R. Administração da Produção e de Operações. 1 ed. Porto Alegre. Bookman. 2009.CORREA, H.L.; CORREA, C.A. Administração da Produção e Operações. 2 ed. São Paulo. Atlas. 2006"
